# Weekly update: insert a new price-report row for "Perejil" (Mercado Mayorista
# Lo Valledor de Santiago) dated 2021-10-07 (serial 44476) just above the
# existing row 305, pushing every subsequent row down by one (305-356 -> 306-357).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 305 - this shifts rows 305:356 down to 306:357
# and Excel carries formatting (e.g. the date style on column D) down from
# the row that used to be above the insertion point.
$ws.Rows.Item(305).EntireRow.Insert()

# Populate the newly inserted row 305 with the new observation.
$ws.Cells.Item(305, 1).Value  = 6
$ws.Cells.Item(305, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(305, 3).Value  = "Metropolitana"
$ws.Cells.Item(305, 4).Value  = 44476
$ws.Cells.Item(305, 5).Value  = 13
$ws.Cells.Item(305, 6).Value  = 100112044
$ws.Cells.Item(305, 7).Value  = "Perejil"
$ws.Cells.Item(305, 8).Value  = "Sin especificar"
$ws.Cells.Item(305, 9).Value  = "Primera"
$ws.Cells.Item(305, 10).Value = 260
$ws.Cells.Item(305, 11).Value = 7000
$ws.Cells.Item(305, 12).Value = 8000
$ws.Cells.Item(305, 13).Value = 7423
$ws.Cells.Item(305, 14).Value = "`$/docena de atados"
$ws.Cells.Item(305, 15).Value = "Región Metropolitana"
$ws.Cells.Item(305, 16).Value = 2474
$ws.Cells.Item(305, 17).Value = 3
$ws.Cells.Item(305, 18).Value = "Hortaliza"
